# The "Restaurant" sheet was missing the "menu" value (column E) for every
# row below the third data row. Fill column E, rows 4-41, with the same
# menu text already used in row 3 ("제육") - effectively a fill-down of E3
# over E4:E41 - and leave the sheet scrolled/selected over that range, the
# way it was left after making the edit interactively.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Restaurant")
[void]$ws.Activate()

$fillValue = $ws.Range("E3").Value()

for ($r = 4; $r -le 41; $r++) {
    $ws.Cells.Item($r, 5).Value = $fillValue
}

[void]$ws.Range("E3:E41").Select()
